$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date strings in column A (rows 3-21): replace "/" separators with "-"
for ($r = 3; $r -le 21; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($val -ne $null) {
        $cell.NumberFormat = "@"
        $cell.Value = $val.ToString().Replace("/", "-")
        $cell.Style = "Normal"
    }
}

# Update the specific attendance counter values that changed for these dates
$ws.Cells.Item(3, 4).Value = 1   # D3
$ws.Cells.Item(3, 7).Value = 1   # G3

$ws.Cells.Item(4, 4).Value = 1   # D4
$ws.Cells.Item(4, 5).Value = 1   # E4
$ws.Cells.Item(4, 8).Value = 0   # H4

$ws.Cells.Item(5, 4).Value = 1   # D5
$ws.Cells.Item(5, 5).Value = 1   # E5
$ws.Cells.Item(5, 8).Value = 0   # H5

$ws.Cells.Item(10, 4).Value = 1  # D10
$ws.Cells.Item(10, 5).Value = 1  # E10
$ws.Cells.Item(10, 8).Value = 0  # H10

$ws.Cells.Item(12, 4).Value = 1  # D12
$ws.Cells.Item(12, 5).Value = 1  # E12
$ws.Cells.Item(12, 8).Value = 0  # H12
